# Update "Name and Address Parsing PPT.pptx":
# On the "Getting the Tokens and Generating the Masks" slide, the body text
#   "Now Suppose the Address Was Split and Tokenized as Follows"
# becomes
#   "Now Suppose the Address and Name Was Split and Tokenized as Follows"
# i.e. the words "and Name " are inserted right after "the Address ".

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -like "*Now Suppose the Address Was Split and Tokenized as Follows*") {
                $targetSlide = $s
                $targetShape = $sh
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$para = $tr.Paragraphs(1, 1)

# "Now Suppose the Address " is 25 characters long; replace the 12
# characters of "the Address " with "the Address and Name " so the
# original run is split around the freshly-typed insertion.
$selection = $para.Characters(13, 12)
$selection.Text = "the Address and Name "
